# Applies the requested change to Plantilla.xlsx:
#  - Insert a new column at AI (shifting FACTURA ACTUAL...UTILIDAD FINAL one column right)
#  - Rename header "PENALIZACIÓN POR DEVOLUCIÓN" (AG8) -> "PENALIZACIÓN POR DEVOLUCIÓN (%)"
#  - Add new header in the freshly inserted column AI8 -> "PENALIZACIÓN POR DEVOLUCIÓN ANTES DE ENTREGA (%)"
#  - Update the active selection to AG8

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before AI (so everything from the old AI onward shifts right by one).
$ws.Range("AI:AI").Insert()

# Update the penalty header text (now a percentage-based penalty).
$ws.Range("AG8").Value = "PENALIZACIÓN POR DEVOLUCIÓN (%)"

# Populate the newly inserted column's header with the new field.
$ws.Range("AI8").Value = "PENALIZACIÓN POR DEVOLUCIÓN ANTES DE ENTREGA (%)"

# Match the style used by neighboring header cells (e.g. AH8 / AJ8).
$ws.Range("AH8").Copy()
$ws.Range("AI8").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = 0

# Re-fit the columns whose header text changed length (AG and the new AI).
$ws.Columns("AG:AG").AutoFit()
$ws.Columns("AI:AI").AutoFit()

# Reflect the new active cell/selection used while editing.
$ws.Application.Goto($ws.Range("AC1"))
$ws.Range("AG8").Select()
